$wb = $excel.ActiveWorkbook
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
